# Update the marksheet's correct/total marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: total marks value 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: corrected/total marks 57 -> 95, and the "x/y" label 55/84 -> 95/140
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
